$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Stall names + prices (rows 2-6 updated in place, rows 7-11 new) ---
$names  = @("Hot Dogs", "Panini", "Tacos", "Pizza", "Ice Cream", "Bubble Tea", "Fried Chips", "Piadina", "Sweets", "Cookies")
$prices = @(5.0, 5.0, 4.0, 10.0, 3.0, 2.0, 3.0, 7.0, 3.0, 3.0)

$ws.Range("B1").Value = "Price"

for ($i = 0; $i -lt $names.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $names[$i]
    $ws.Cells.Item($row, 2).Value = $prices[$i]
}

# New rows 7-11 need the same look as the existing data rows (4-6):
# copy that formatting down before touching the now-obsolete X/Y/Type columns.
$ws.Range("A6:B6").Copy()
$ws.Range("A7:B11").PasteSpecial(-4122)

# --- Drop the old X / Y / Type columns -------------------------------
# Column D (old "Price") is removed outright.
$ws.Range("D1:D6").Clear()
# Column C (old "Y") keeps its row-1..3 formatting but loses its data,
# and is dropped completely for rows 4-6 (no data ever lived beyond row 3).
$ws.Range("C1:C3").ClearContents()
$ws.Range("C4:C6").Clear()
# Column E (old "Type") keeps its formatting everywhere but loses its data.
$ws.Range("E1:E6").ClearContents()
